$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range of the worksheet to find the last row with data.
$lastRow = $ws.Cells.SpecialCells(11).Row  # xlCellTypeLastCell = 11

# Data rows start at row 2 (row 1 is the header row).
$startRow = 2

# Column C ("Förändrad") holds a date serial number that needs to be
# incremented by one day (45177 -> 45178) for every data row.
$range = $ws.Range("C$startRow`:C$lastRow")
$range.Value = 45178
